$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 37 (pushes existing rows 37..92 down to 38..93)
$ws.Rows.Item(37).Insert()

# Populate the newly inserted row 37 with the new weekly price observation
$ws.Cells.Item(37, 1).Value = 4
$ws.Cells.Item(37, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(37, 3).Value = "Los Lagos"
$ws.Cells.Item(37, 4).Value = 44413
$ws.Cells.Item(37, 5).Value = 10
$ws.Cells.Item(37, 6).Value = 100112039
$ws.Cells.Item(37, 7).Value = "Ciboulette"
$ws.Cells.Item(37, 8).Value = "Sin especificar"
$ws.Cells.Item(37, 9).Value = "Primera"
$ws.Cells.Item(37, 10).Value = 120
$ws.Cells.Item(37, 11).Value = 4500
$ws.Cells.Item(37, 12).Value = 4500
$ws.Cells.Item(37, 13).Value = 4500
$ws.Cells.Item(37, 14).Value = "$/docena de atados"
$ws.Cells.Item(37, 15).Value = "Región Metropolitana"
$ws.Cells.Item(37, 16).Value = 1500
$ws.Cells.Item(37, 17).Value = 3
$ws.Cells.Item(37, 18).Value = "Hortaliza"
